# Update Resumo de Inscricoes - adjust inscription counts for several vagas
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Inscricoes")

$ws.Range("E7").Value = 15

$ws.Range("E10").Value = 210
$ws.Range("F10").Value = 90
$ws.Range("H10").Value = 90

$ws.Range("E11").Value = 156

$ws.Range("E12").Value = 228
$ws.Range("F12").Value = 115
$ws.Range("H12").Value = 115

$ws.Range("E18").Value = 32

$ws.Range("E21").Value = 73
$ws.Range("F21").Value = 37
$ws.Range("H21").Value = 37

$ws.Range("E22").Value = 92
$ws.Range("F22").Value = 38
$ws.Range("H22").Value = 38

$ws.Range("F23").Value = 45
$ws.Range("H23").Value = 45

$ws.Range("E24").Value = 106
$ws.Range("F24").Value = 50
$ws.Range("H24").Value = 50

$ws.Range("E25").Value = 103
$ws.Range("F25").Value = 37
$ws.Range("H25").Value = 37

$ws.Range("E26").Value = 62
$ws.Range("F26").Value = 28
$ws.Range("H26").Value = 28

$ws.Range("E27").Value = 152

$ws.Range("E28").Value = 93
$ws.Range("F28").Value = 26
$ws.Range("H28").Value = 26

$ws.Range("E29").Value = 97
$ws.Range("F29").Value = 54
$ws.Range("H29").Value = 54

$ws.Range("E30").Value = 110

$ws.Range("E32").Value = 104

$ws.Range("E33").Value = 138
$ws.Range("F33").Value = 63
$ws.Range("H33").Value = 63

$ws.Range("E34").Value = 109

$ws.Range("E35").Value = 70

$ws.Range("E36").Value = 32

$ws.Range("E38").Value = 49
$ws.Range("F38").Value = 26
$ws.Range("H38").Value = 26

$ws.Range("E41").Value = 186

$ws.Range("E42").Value = 159
$ws.Range("F42").Value = 74
$ws.Range("H42").Value = 74

$ws.Range("E45").Value = 60

$ws.Range("E46").Value = 128

$ws.Range("E47").Value = 216
$ws.Range("F47").Value = 87
$ws.Range("H47").Value = 87

$ws.Range("E48").Value = 107

$ws.Range("E49").Value = 115

$ws.Range("E50").Value = 99
$ws.Range("F50").Value = 35
$ws.Range("H50").Value = 35

$ws.Range("E51").Value = 104
$ws.Range("F51").Value = 38
$ws.Range("H51").Value = 38
